$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date column (column B) for data rows 2 through 175
# from the old date (45139/45140, i.e. 2023-08-01/02) to the new
# date 45693 (2025-02-05).
$ws.Range("B2:B175").Value = 45693

# Update the frozen-pane top-left cell and selection to reflect the
# new active cell/selection location.
$ws.Range("B2:B175").Select()
